$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.ClearFormats()
}

Set-TextValue $ws.Range("D2") '27.205.54'
Set-TextValue $ws.Range("E2") '  -2.34%  '
Set-TextValue $ws.Range("D3") '1.872.72'
Set-TextValue $ws.Range("E3") '  -1.77%  '
Set-TextValue $ws.Range("D4") '1.002'
Set-TextValue $ws.Range("E4") '  -0.27%  '
Set-TextValue $ws.Range("D5") '307.38'
Set-TextValue $ws.Range("E5") '  -2.03%  '
Set-TextValue $ws.Range("D6") '1.000'
Set-TextValue $ws.Range("E6") '  -0.29%  '
Set-TextValue $ws.Range("D7") '0.5120'
Set-TextValue $ws.Range("E7") '  +1.92%  '
Set-TextValue $ws.Range("D8") '0.3760'
Set-TextValue $ws.Range("E8") '  -1.50%  '
Set-TextValue $ws.Range("D9") '0.07169'
Set-TextValue $ws.Range("E9") '  -1.39%  '
Set-TextValue $ws.Range("D10") '0.8896'
Set-TextValue $ws.Range("E10") '  -2.04%  '
Set-TextValue $ws.Range("E11") '  -0.49%  '
Set-TextValue $ws.Range("E12") '  -1.06%  '
Set-TextValue $ws.Range("D13") '1.851.25'
Set-TextValue $ws.Range("E13") '  -3.30%  '
Set-TextValue $ws.Range("D14") '5.344'
Set-TextValue $ws.Range("E14") '  -2.42%  '
Set-TextValue $ws.Range("D15") '89.55'
Set-TextValue $ws.Range("E15") '  -2.44%  '
Set-TextValue $ws.Range("D16") '1.001'
Set-TextValue $ws.Range("E16") '  -0.37%  '
Set-TextValue $ws.Range("D17") '0.000008558'
Set-TextValue $ws.Range("E17") '  -1.73%  '
Set-TextValue $ws.Range("E18") '  -2.64%  '
Set-TextValue $ws.Range("D19") '0.9999'
Set-TextValue $ws.Range("E19") '  -0.36%  '
Set-TextValue $ws.Range("D20") '27.253.63'
Set-TextValue $ws.Range("E20") '  -2.30%  '
Set-TextValue $ws.Range("D21") '5.076'
Set-TextValue $ws.Range("E21") '  -1.75%  '
Set-TextValue $ws.Range("D22") '2.094.28'
Set-TextValue $ws.Range("E22") '  -3.30%  '
Set-TextValue $ws.Range("D23") '10.64'
Set-TextValue $ws.Range("E23") '  -1.61%  '
Set-TextValue $ws.Range("D24") '6.503'
Set-TextValue $ws.Range("E24") '  -1.13%  '
Set-TextValue $ws.Range("D25") '151.04'
Set-TextValue $ws.Range("E25") '  -1.90%  '
Set-TextValue $ws.Range("D26") '1.847'
Set-TextValue $ws.Range("E26") '  -1.53%  '
Set-TextValue $ws.Range("D27") '18.02'
Set-TextValue $ws.Range("E27") '  -1.98%  '
Set-TextValue $ws.Range("D28") '2.134'
Set-TextValue $ws.Range("E28") '  -3.62%  '
Set-TextValue $ws.Range("D29") '112.67'
Set-TextValue $ws.Range("E29") '  -2.19%  '
Set-TextValue $ws.Range("D30") '4.760'
Set-TextValue $ws.Range("E30") '  -2.91%  '
Set-TextValue $ws.Range("D31") '4.708'
Set-TextValue $ws.Range("E31") '  +1.30%  '
Set-TextValue $ws.Range("D32") '0.08998'
Set-TextValue $ws.Range("E32") '  -0.05%  '
Set-TextValue $ws.Range("D33") '0.05166'
Set-TextValue $ws.Range("E33") '  -1.51%  '
Set-TextValue $ws.Range("D34") '3.107'
Set-TextValue $ws.Range("E34") '  -3.35%  '
Set-TextValue $ws.Range("D35") '0.7551'
Set-TextValue $ws.Range("E35") '  -0.74%  '
Set-TextValue $ws.Range("D36") '1.174'
Set-TextValue $ws.Range("E36") '  -4.08%  '
Set-TextValue $ws.Range("D37") '0.02045'
Set-TextValue $ws.Range("E37") '  -0.75%  '
Set-TextValue $ws.Range("E38") '  +0.91%  '
Set-TextValue $ws.Range("D39") '3.036'
Set-TextValue $ws.Range("E39") '  +0.42%  '
Set-TextValue $ws.Range("D41") '0.5360'
Set-TextValue $ws.Range("E41") '  -3.44%  '
Set-TextValue $ws.Range("D42") '6.668'
Set-TextValue $ws.Range("E42") '  -3.16%  '
Set-TextValue $ws.Range("D43") '114.67'
Set-TextValue $ws.Range("E43") '  +3.22%  '
Set-TextValue $ws.Range("D44") '8.566'
Set-TextValue $ws.Range("E44") '  +1.10%  '
Set-TextValue $ws.Range("D45") '0.1485'
Set-TextValue $ws.Range("E45") '  -1.67%  '
Set-TextValue $ws.Range("D46") '0.4683'
Set-TextValue $ws.Range("E46") '  -2.83%  '
Set-TextValue $ws.Range("D47") '1.000'
Set-TextValue $ws.Range("E47") '  -0.33%  '
Set-TextValue $ws.Range("D48") '10.07'
Set-TextValue $ws.Range("E48") '  -4.72%  '
Set-TextValue $ws.Range("D49") '1.574'
Set-TextValue $ws.Range("E49") '  -3.13%  '
Set-TextValue $ws.Range("D50") '65.24'
Set-TextValue $ws.Range("E50") '  -3.12%  '
Set-TextValue $ws.Range("D51") '36.64'
Set-TextValue $ws.Range("E51") '  -0.81%  '
